# Update TPM-derived LR-pair metrics (columns E:T) for data rows 2-10
# on Sheet1, per the "update scripts wuth new tpm" recompute.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$firstDataRow = 2
$firstCol = 5   # column E
$lastCol  = 20  # column T

$newValues = @(
    @(3, 1, 2.311181333333333, 6.933544, 0.336172840858964, 0.336172840858964, 2, 0.6666666666666666, 0.06247066666666667, 0.187412, 0.009697282539039283, 0.009697282539039282, 0.1443810386808889, 1.299429348128, 0.003259963019760863, 0.003259963019760863),
    @(3, 1, 2.311181333333333, 6.933544, 0.336172840858964, 0.336172840858964, 3, 1, 3.384345, 10.153035, 0.5253497589468907, 0.5253497589468907, 7.821834989559999, 70.39651490604, 0.1766083209097482, 0.1766083209097482),
    @(3, 1, 2.311181333333333, 6.933544, 0.336172840858964, 0.336172840858964, 3, 1, 2.995264, 8.985792, 0.46495295851407, 0.4649529585140699, 6.922598245205333, 62.303384206848, 0.1563045569294549, 0.1563045569294549),
    @(3, 1, 4.100181333333333, 12.300544, 0.5963918049111226, 0.5963918049111226, 2, 0.6666666666666666, 0.06247066666666667, 0.187412, 0.009697282539039283, 0.009697282539039282, 0.2561410613475555, 2.305269552128, 0.005783379836190752, 0.005783379836190751),
    @(3, 1, 4.100181333333333, 12.300544, 0.5963918049111226, 0.5963918049111226, 3, 1, 3.384345, 10.153035, 0.5253497589468907, 0.5253497589468907, 13.87642819456, 124.88785375104, 0.3133142909479594, 0.3133142909479594),
    @(3, 1, 4.100181333333333, 12.300544, 0.5963918049111226, 0.5963918049111226, 3, 1, 2.995264, 8.985792, 0.46495295851407, 0.4649529585140699, 12.28112554120533, 110.530129870848, 0.2772941341269725, 0.2772941341269725),
    @(3, 1, 0.4636166666666667, 1.39085, 0.06743535422991333, 0.06743535422991333, 2, 0.6666666666666666, 0.06247066666666667, 0.187412, 0.009697282539039283, 0.009697282539039282, 0.02896244224444445, 0.2606619802, 0.0006539396830876674, 0.0006539396830876672),
    @(3, 1, 0.4636166666666667, 1.39085, 0.06743535422991333, 0.06743535422991333, 3, 1, 3.384345, 10.153035, 0.5253497589468907, 0.5253497589468907, 1.56903874775, 14.12134872975, 0.03542714708918315, 0.03542714708918315),
    @(3, 1, 0.4636166666666667, 1.39085, 0.06743535422991333, 0.06743535422991333, 3, 1, 2.995264, 8.985792, 0.46495295851407, 0.4649529585140699, 1.388654311466667, 12.4978888032, 0.03135426745764251, 0.0313542674576425)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowNum = $firstDataRow + $i
    $rowVals = $newValues[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($rowNum, $firstCol + $c).Value = $rowVals[$c]
    }
}
